$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels (A1:P1)
$headerRow = @("Language", "German", "Spanish", "Slovak", "Norwegian", "Greek", "Chinese", "Vietnamese", "Indonesian", "Finnish", "Basque", "Korean", "Japanese", "Turkish", "Arabic", "Hebrew")
for ($c = 0; $c -lt $headerRow.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headerRow[$c]
}

# Rows 2-17: language label in col A, then 15 numeric similarity values
$matrix = @(
    @("German", 0.9824459770114943, 0.8229394636015326, 0.7516291187739463, 0.8470927203065134, 0.7861885057471264, 0.5591141762452108, 0.574663601532567, 0.7238620689655172, 0.667656704980843, 0.6131157088122605, 0.5895019157088123, 0.5680674329501916, 0.625072796934866, 0.6401900383141762, 0.7131310344827586),
    @("Spanish", 0.844756728092271, 0.9882388591125169, 0.777674665252552, 0.8765174901990492, 0.7590954716766728, 0.5635499327664248, 0.5846858961004526, 0.8137535273953145, 0.6331319482585557, 0.5701785951023655, 0.5201984810893733, 0.5320732561883298, 0.5775458798128823, 0.6349311566068825, 0.7232628169921024),
    @("Slovak", 0.829290758366595, 0.8850936444580902, 0.9656892846177464, 0.8583819465766043, 0.8271415412956709, 0.6410807491556647, 0.6007829290758366, 0.7793214614676083, 0.77655818237642, 0.723672090881179, 0.606002456248081, 0.6220448265274793, 0.6947344181762358, 0.7480810561866749, 0.7457783236106846),
    @("Norwegian", 0.8738570379763733, 0.8544683975171862, 0.7756457318294067, 0.9789094306881132, 0.7917639991990923, 0.6221384235466862, 0.5662751117933658, 0.8012414069278516, 0.7456116932523527, 0.6140959754388307, 0.5691116598811987, 0.5502903290395782, 0.5960421811386237, 0.6836748314756724, 0.7661349529466729),
    @("Greek", 0.8044522445963421, 0.7745242933678183, 0.73766857565121, 0.7579900240162571, 0.9757990024016258, 0.5254018104563088, 0.5047108812118973, 0.6810456308886015, 0.6591538887862554, 0.5773138740070202, 0.5190282652872714, 0.5503417698134122, 0.6302420099759838, 0.6095510807315722, 0.6831701459449473),
    @("Chinese", 0.5141525141525142, 0.567016317016317, 0.4023476523476524, 0.6224608724608724, 0.5213120213120214, 0.9492174492174492, 0.5404595404595405, 0.5937395937395937, 0.5826673326673326, 0.5621878121878122, 0.5382950382950383, 0.5938228438228438, 0.5105727605727606, 0.4392274392274392, 0.5322177822177823),
    @("Vietnamese", 0.5406942701798411, 0.5872020075282308, 0.5648682559598495, 0.5613550815558344, 0.5371810957758261, 0.5622751986616479, 0.8945211208699289, 0.6271016311166876, 0.550564617314931, 0.5733166039314095, 0.5038895859473024, 0.5169385194479298, 0.5553324968632372, 0.5375993308239231, 0.5989125888749477),
    @("Thai", 0.3516709972224711, 0.4116566615894633, 0.444270226682197, 0.3193262252486336, 0.4361616342621629, 0.463981722067915, 0.5140220410357494, 0.3885852522175433, 0.4663112624316818, 0.4727622972851895, 0.3951258847773497, 0.4317713466535257, 0.4454349968640803, 0.4371920078845982, 0.5891497177672251),
    @("Indonesian", 0.7381154499151104, 0.8202886247877759, 0.750169779286927, 0.8225806451612904, 0.6646010186757215, 0.6289473684210526, 0.580730050933786, 0.9275891341256367, 0.7298811544991511, 0.6679966044142615, 0.5992359932088285, 0.5795415959252971, 0.6387945670628183, 0.61553480475382, 0.7343803056027165),
    @("Finnish", 0.7401306447032093, 0.8178074410678784, 0.739136608917921, 0.7874183470604942, 0.7884597178831771, 0.6336741456025751, 0.5858657578339487, 0.7567452428287418, 0.9651614124775159, 0.7496923222569346, 0.6361355675470983, 0.6616964877402253, 0.7340244248792956, 0.636324907696677, 0.7894537536684654),
    @("Basque", 0.6565192418150488, 0.6682940838598507, 0.6449905637154345, 0.6432263887749241, 0.6327644211044555, 0.5720439812915402, 0.5607614671371133, 0.6954951998030688, 0.7055468942315583, 0.9481414622138344, 0.6275539509313203, 0.6913104127348815, 0.6776483137769754, 0.5816443751538525, 0.7034545006974645),
    @("Korean", 0.5504124656278644, 0.565853486568427, 0.5163223577522386, 0.556581823309596, 0.5499189170133258, 0.467531551857858, 0.484206444334767, 0.568638510893323, 0.5782979623492914, 0.571494042163153, 0.9524430656419657, 0.4936896284284002, 0.5566170767820631, 0.4608686455615878, 0.5665585560177677),
    @("Japanese", 0.3458646616541353, 0.4875709682369188, 0.3107257940770293, 0.3551480742673009, 0.3672702163572196, 0.4797452815712751, 0.4478287555623753, 0.4375479515114317, 0.4686205309191346, 0.5049102347705999, 0.4796685591529845, 0.9703851465398189, 0.518796992481203, 0.3775510204081632, 0.4825840110480282),
    @("Turkish", 0.605129705480788, 0.655061439438268, 0.6150770431051297, 0.6090306221962161, 0.6302906182952994, 0.6277550224302711, 0.6203432806709577, 0.6762239126194656, 0.6914374878096352, 0.6979715233079774, 0.5769455822118198, 0.6303881412131851, 0.941876340940121, 0.5473961380924517, 0.6609128145114102),
    @("Arabic", 0.6261357378854625, 0.6510187224669604, 0.6603111233480177, 0.6542194383259912, 0.662961178414097, 0.4665129405286343, 0.5424008810572687, 0.6526362885462555, 0.5567524779735683, 0.5299765969162996, 0.4578400330396476, 0.4740501101321586, 0.5228180066079295, 0.9678207599118943, 0.7954983480176211),
    @("Hebrew", 0.559291801545881, 0.6430600515293652, 0.5401995111316641, 0.5496465614058268, 0.5993922177445993, 0.4410385148972716, 0.469709982162912, 0.5704564973244368, 0.4909823611019357, 0.4625751469908172, 0.4131598070951972, 0.4547136156437868, 0.4656140582678206, 0.6397568870978397, 0.9672986721278986)
)

for ($r = 0; $r -lt $matrix.Length; $r++) {
    $row = $matrix[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Update the saved selection/active cell to match the authored edit
$ws.Range("D23").Select()